$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.000", "30.515.91") that must
# stay literal text rather than being auto-coerced to a number by Excel.
# Pre-format the whole data range as Text so assignment below keeps the string.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.515.91'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.852.86'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '233.92'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.4707'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').Value = '0.2746'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.06366'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '17.70'
$ws.Range('E10').Value = '  +8.41%  '
$ws.Range('D11').Value = '1.877.43'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '0.07417'
$ws.Range('D13').Value = '5.046'
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = '84.70'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = '0.6269'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '30.471.19'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '242.34'
$ws.Range('E17').Value = '  +5.39%  '
$ws.Range('D18').Value = '0.9997'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Value = '12.71'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = '0.000007357'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = '0.9999'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').Value = '4.961'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').Value = '5.992'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '9.247'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').Value = '162.10'
$ws.Range('E25').Value = '  -2.57%  '
$ws.Range('D26').Value = '18.03'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').Value = '1.891'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').Value = '0.1019'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = '1.366'
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').Value = '4.025'
$ws.Range('E30').Value = '  -3.12%  '
$ws.Range('D31').Value = '3.852'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').Value = '0.04888'
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('D33').Value = '1.139'
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('D34').Value = '0.7081'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '0.01902'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '2.688'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('D38').Value = '0.8743'
$ws.Range('E38').Value = '  -4.72%  '
$ws.Range('D39').Value = '1.978'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '105.40'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '0.4084'
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').Value = '5.509'
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('D44').Value = '7.226'
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('D45').Value = '62.60'
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('D46').Value = '0.1208'
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '8.581'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '33.35'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.379'
$ws.Range('E49').Value = '  -2.04%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05545'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '0.3689'
$ws.Range('E51').Value = '  -0.61%  '
